# "Update results for Steel"
# Sheet1 layout: columns B/C/D = Iron & steel / Chemicals / Non-metallic minerals
# rows 3/5/6/8 = Hydrogen / Ammonia / Biomass / Other
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iron & steel - Hydrogen result updated
$ws.Range("B3").Value = 156.5814015774635

# Minor recalculation refresh on the other industries' results
$ws.Range("C5").Value = 2977.411704809431
$ws.Range("D6").Value = 910.7638679465082
$ws.Range("D8").Value = 393.9265512588963
